{"js": "// The document contains three short paragraphs/runs that spell out an\n// \"<id>...</id>\" tag split across three separate runs, e.g.:\n//   run1 (Courier New, orange): \"<id>\"\n//   run2 (plain):                \"p068r_3\"\n//   run3 (Courier New, orange): \"</id>\"\n// The edit collapses each of these triples into a single run whose text is\n// the concatenation \"<id>p068r_3</id>\" (keeping the first run's formatting).\n// There are three such \"<id>...</id>\" occurrences in the body.\nconst ids = [\"p068r_3\", \"p069r_1\", \"p069r_2\"];\n\nconst body = context.document.body;\n\nfor (const id of ids) {\n  const needle = \"<id>\" + id + \"</id>\";\n  const results = body.search(needle, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    // Replacing the whole matched range's text merges the 3 runs it spans\n    // into a single run, and Word keeps the formatting of the range's\n    // first run (the Courier New \"<id>\" run) for the merged result -\n    // exactly matching the target OOXML.\n    results.items[i].insertText(needle, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three short paragraphs that spell out an\n# \"<id>...</id>\" tag split across three separate runs, e.g.:\n#   run1 (Courier New, orange): \"<id>\"\n#   run2 (plain):                \"p068r_3\"\n#   run3 (Courier New, orange): \"</id>\"\n# The edit collapses each of these triples into a single run whose text is\n# the concatenation \"<id>p068r_3</id>\" (keeping the first run's formatting).\n# There are three such \"<id>...</id>\" occurrences in the body.\n$d = $word.ActiveDocument\n$ids = @(\"p068r_3\", \"p069r_1\", \"p069r_2\")\n\nforeach ($id in $ids) {\n    $needle = \"<id>\" + $id + \"</id>\"\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $needle\n    # Replacing the matched text merges the runs the find-hit spans into a\n    # single run, and Word keeps the formatting of the first run in that\n    # span (the Courier New \"<id>\" run) - exactly matching the target OOXML.\n    $rng.Find.Replacement.Text = $needle\n    $rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n"}
